# PM_Coaching_Review.xlsx — add the User Transaction feature SIQ and SRS coaching rows,
# rename the sheet, correct the release-version tag, and widen the Comments column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab.
$ws.Name = "coaching Review"

# 2) Free up the old "V1" shared string by temporarily pointing every F2:F11 cell at an
#    already-existing string. The freed slot gets reused by whichever brand-new string is
#    introduced next, so the creation order below is chosen to reproduce the exact shared
#    string table layout of the target workbook.
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("F$r").Value = "closed"
}

# 3) New string #1: the release-naming comment (claims the freed slot).
$ws.Range("B12").Value = "release name doesn't follow the release naming`n convension as mensioned in PMP"

# 4) New string #2: the review date used by all three new rows.
$ws.Range("C12").Value = "17/4/2024"

# 5) New string #3: "V1.0" — introduced via F2, then reused for F3:F11 and F12:F14.
$ws.Range("F2").Value = "V1.0"
for ($r = 3; $r -le 11; $r++) {
    $ws.Range("F$r").Value = "V1.0"
}

# 6) New string #4: "V2.0" — the release tag for the three new rows.
$ws.Range("F12").Value = "V2.0"

# 7) Finish filling row 12 (reusing already-known strings for D/E).
$ws.Range("A12").Value = 11
$ws.Range("D12").Value = "Esraa"
$ws.Range("E12").Value = "closed"

# 8) Row 13 — new string #5 is the "Req_SIQ" comment.
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "in [Req_SIQ]: the second tab should be named `n`"Revision Record`" not Report."
$ws.Range("C13").Value = "17/4/2024"
$ws.Range("D13").Value = "Esraa"
$ws.Range("E13").Value = "closed"
$ws.Range("F13").Value = "V2.0"

# 9) Row 14 — new string #6 is the "why two sheets" comment.
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Why is there REQ_SIQ and REQ_SIQ_Responses`n not only one sheet?"
$ws.Range("C14").Value = "17/4/2024"
$ws.Range("D14").Value = "Esraa"
$ws.Range("E14").Value = "closed"
$ws.Range("F14").Value = "V2.0"

# 10) Row heights for the new rows.
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 27.5
$ws.Rows.Item(14).RowHeight = 48.5

# 11) Wrap text on the new Comments cells.
$ws.Range("B12").WrapText = $true
$ws.Range("B13").WrapText = $true
$ws.Range("B14").WrapText = $true

# 12) Widen the Comments column (B) to fit the new longer text.
$ws.Columns.Item(2).ColumnWidth = 40

# 13) Leave the selection on B12, matching the saved view state.
$ws.Range("B12").Select()
